$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Complete the existing last row (181): it already has A181 (Serie) and
#     G181 (Tasa de politica monetaria); fill in the other columns.
$ws.Range("B181").Value = 80000
$ws.Range("C181").Value = 1.5
$ws.Range("D181").Value = 1.5
$ws.Range("E181").Value = 1.5
$ws.Range("F181").Value = 5
$ws.Range("G181").Value = 1.5

# --- New daily rows 182-190 (21-09-2021 .. 01-10-2021) ---
# Columns: A=Serie(date) B=Monto C=Tasa promedio D=Tasa maxima E=Tasa minima
#          F=Participantes G=Tasa politica monetaria
$dates = @(
    "21-09-2021",
    "22-09-2021",
    "23-09-2021",
    "24-09-2021",
    "27-09-2021",
    "28-09-2021",
    "29-09-2021",
    "30-09-2021",
    "01-10-2021"
)

$rowsData = @(
    @{ B = 110000; C = 1.5; D = 1.5; E = 1.5; F = 4; G = 1.5 },
    @{ B = 40000;  C = 1.5; D = 1.5; E = 1.5; F = 3; G = 1.5 },
    @{ B = 180000; C = 1.5; D = 1.5; E = 1.5; F = 4; G = 1.5 },
    @{ B = 70000;  C = 1.5; D = 1.5; E = 1.5; F = 4; G = 1.5 },
    @{ B = 85000;  C = 1.5; D = 1.5; E = 1.5; F = 3; G = 1.5 },
    @{ B = 85000;  C = 1.5; D = 1.5; E = 1.5; F = 3; G = 1.5 },
    @{ B = 80000;  C = 1.5; D = 1.5; E = 1.5; F = 3; G = 1.5 },
    @{ B = 20000;  F = 2;   G = 1.5 },
    @{ G = 1.5 }
)

$startRow = 182
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $dateText = $dates[$i]

    # Some of these day strings (e.g. "01-10-2021") look like an unambiguous
    # MM-DD-YYYY date to Excel's smart cell-entry logic, so a plain
    # Range.Value assignment would silently convert it into a date serial
    # number instead of keeping it as the literal text used everywhere else
    # in column A. Route it through a text formula + paste-values round
    # trip so the literal string lands in the cell no matter what it looks
    # like, without leaving any stray number-format behind.
    $ws.Range("ZZ1").Formula = '="' + $dateText + '"'
    $ws.Range("ZZ1").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4163) | Out-Null
    $ws.Range("ZZ1").Clear() | Out-Null

    $data = $rowsData[$i]
    if ($data.ContainsKey("B")) { $ws.Range("B$r").Value = $data.B }
    if ($data.ContainsKey("C")) { $ws.Range("C$r").Value = $data.C }
    if ($data.ContainsKey("D")) { $ws.Range("D$r").Value = $data.D }
    if ($data.ContainsKey("E")) { $ws.Range("E$r").Value = $data.E }
    if ($data.ContainsKey("F")) { $ws.Range("F$r").Value = $data.F }
    if ($data.ContainsKey("G")) { $ws.Range("G$r").Value = $data.G }
}

$excel.CutCopyMode = $false
